$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.66"
$ws.Range("E2").Value = "'-2.99%"
$ws.Range("D3").Value = "'37.24"
$ws.Range("E3").Value = "'-6.64%"
$ws.Range("D4").Value = "'5.112"
$ws.Range("E4").Value = "'-0.40%"
$ws.Range("D5").Value = "'0.07819"
$ws.Range("E5").Value = "'-4.82%"
$ws.Range("D6").Value = "'1.957"
$ws.Range("E6").Value = "'-4.57%"
$ws.Range("D7").Value = "'4.385"
$ws.Range("E7").Value = "'2.13%"
$ws.Range("D8").Value = "'8.278"
$ws.Range("E8").Value = "'-0.49%"
$ws.Range("D9").Value = "'3.031"
$ws.Range("E9").Value = "'-8.40%"
$ws.Range("D10").Value = "'0.9250"
$ws.Range("E10").Value = "'-0.99%"
$ws.Range("D11").Value = "'0.1332"
$ws.Range("E11").Value = "'-2.36%"
$ws.Range("D12").Value = "'0.1945"
$ws.Range("E12").Value = "'-1.84%"
$ws.Range("D13").Value = "'0.09004"
$ws.Range("E13").Value = "'-0.48%"
$ws.Range("D14").Value = "'0.03442"
$ws.Range("D15").Value = "'0.09709"
$ws.Range("E15").Value = "'-0.93%"
$ws.Range("D16").Value = "'0.001391"
$ws.Range("E16").Value = "'-1.03%"
$ws.Range("D17").Value = "'0.005932"
$ws.Range("E17").Value = "'-6.04%"
$ws.Range("D18").Value = "'3.593"
$ws.Range("E18").Value = "'-2.37%"
$ws.Range("D19").Value = "'0.3424"
$ws.Range("E19").Value = "'-1.46%"
$ws.Range("E20").Value = "'0.10%"
$ws.Range("D21").Value = "'5.005"
$ws.Range("E21").Value = "'2.32%"
$ws.Range("D22").Value = "'0.2491"
$ws.Range("E22").Value = "'1.88%"
$ws.Range("D23").Value = "'0.02106"
$ws.Range("E23").Value = "'5,178.69%"
$ws.Range("D24").Value = "'0.04333"
$ws.Range("E24").Value = "'0.08%"
$ws.Range("D25").Value = "'0.001217"
$ws.Range("E25").Value = "'-0.70%"
$ws.Range("D26").Value = "'0.004524"
$ws.Range("E26").Value = "'-5.21%"
$ws.Range("D27").Value = "'0.0001351"
$ws.Range("E27").Value = "'4.10%"
$ws.Range("D39").Value = "'0.02286"
$ws.Range("E39").Value = "'2.72%"
$ws.Range("D40").Value = "'0.05025"
$ws.Range("E40").Value = "'-3.76%"
$ws.Range("E41").Value = "'0.51%"
$ws.Range("D42").Value = "'0.009808"
$ws.Range("E42").Value = "'0.77%"
$ws.Range("D43").Value = "'0.1350"
$ws.Range("E43").Value = "'-2.60%"
$ws.Range("D44").Value = "'0.002062"
$ws.Range("E44").Value = "'-1.36%"
$ws.Range("D45").Value = "'0.008447"
$ws.Range("E45").Value = "'-7.99%"
$ws.Range("D46").Value = "'0.00006778"
$ws.Range("E46").Value = "'3.47%"
$ws.Range("E47").Value = "'0.35%"
$ws.Range("D48").Value = "'0.003017"
$ws.Range("E48").Value = "'1.15%"
$ws.Range("E49").Value = "'-22.79%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.35%"
$ws.Range("E51").Value = "'0.35%"
